# Commit: "Fruta / hortaliza, semanal"
# Insert a new weekly data row at row 8 (pushing existing rows 8..77 down to 9..78)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8; this shifts rows 8-77 down to 9-78
# and Excel automatically extends the used range / dimension (A1:R77 -> A1:R78).
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44537
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = "Arveja Verde"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 18000
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Región de La Araucanía"
$ws.Range("P8").Value = 720
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
